$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.336.98"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").Value = "2.935.06"
$ws.Range("E3").Value = "  +0.75%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "358.44"
$ws.Range("E5").Value = "  +1.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.22"
$ws.Range("E6").Value = "  -2.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("E7").Value = "  +2.11%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.635"
$ws.Range("E9").Value = "  +0.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.11"
$ws.Range("E10").Value = "  -2.40%  "

$ws.Range("E11").Value = "  +1.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0872"
$ws.Range("E12").Value = "  +0.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.66"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.78"
$ws.Range("E14").Value = "  -0.54%  "

$ws.Range("D15").Value = "3.390.22"
$ws.Range("E15").Value = "  +0.77%  "

$ws.Range("D16").Value = "2.923.54"
$ws.Range("E16").Value = "  +0.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.989"
$ws.Range("E17").Value = "  -1.41%  "

$ws.Range("D18").Value = "52.280.16"
$ws.Range("E18").Value = "  -0.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.55"
$ws.Range("E19").Value = "  +6.96%  "

$ws.Range("E20").Value = "  -0.70%  "

$ws.Range("E21").Value = "  -2.22%  "

$ws.Range("D22").Value = "0.0₃0985"
$ws.Range("E22").Value = "  +0.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.65"
$ws.Range("E23").Value = "  -0.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.83"
$ws.Range("E24").Value = "  -0.56%  "

$ws.Range("E25").Value = "  +1.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.185"
$ws.Range("E26").Value = "  +6.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.07"
$ws.Range("E27").Value = "  +0.77%  "

$ws.Range("E28").Value = "  +15.13%  "

$ws.Range("E29").Value = "  +0.20%  "

$ws.Range("E30").Value = "  +8.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.52"
$ws.Range("E31").Value = "  -1.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.50"
$ws.Range("E32").Value = "  -1.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.21"
$ws.Range("E33").Value = "  -2.66%  "

$ws.Range("E34").Value = "  -1.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "52.31"
$ws.Range("E35").Value = "  -1.89%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0444"
$ws.Range("E36").Value = "  -1.83%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.21"
$ws.Range("E38").Value = "  -3.93%  "

$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.37"
$ws.Range("E39").Value = "  -2.81%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.81"
$ws.Range("E40").Value = "  -1.70%  "

$ws.Range("E41").Value = "  -3.52%  "

$ws.Range("E42").Value = "  +2.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.12"
$ws.Range("E43").Value = "  -1.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "120.16"
$ws.Range("E44").Value = "  -0.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.18"
$ws.Range("E45").Value = "  -1.04%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.47"
$ws.Range("E46").Value = "  -2.49%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.48"
$ws.Range("E47").Value = "  -5.27%  "

$ws.Range("D48").Value = "2.132.09"
$ws.Range("E48").Value = "  -3.20%  "

$ws.Range("E49").Value = "  -5.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0350"
$ws.Range("E50").Value = "  +2.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.934"
$ws.Range("E51").Value = "  -4.19%  "
